$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new question row (row 9): a FILL-in-the-blank question about India's capital
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "FILL"
$ws.Cells.Item(9, 3).Value = "INDIAS CAPITAl?"

# The existing data rows (2-8) lose their old default cell styling
$ws.Range("A2:H8").ClearFormats()

# Re-enter the header so it ends up as the freshest shared string / keeps same text
$ws.Range("H1").Value = "correct_asnwer"

# Resize columns to fit the refreshed content (closest values this host's
# character-width quantization can land on, targeting 84.21875 / 21.6640625 /
# 15.6640625 / 27.33203125 respectively)
$ws.Range("C1").ColumnWidth = 83.255
$ws.Range("D1").ColumnWidth = 20.754062500000046
$ws.Range("E1").ColumnWidth = 14.754062499999995
$ws.Range("F1").ColumnWidth = 26.422031250000046

# Move the selection/cursor to H1, matching where editing finished
$ws.Range("H1").Select()
